# Apply the edit described by the commit:
#  - Reword the "Surface tension coefficient (sigma)" header (column D) to add a
#    stability-issue caveat, and widen the column to fit the new text.
#  - Reset the sigma values (D3, D4) to 0 and change their number format from the
#    old custom decimal formats to a plain integer format ("0").
#  - Fix the x size / y size (N4, O4) values in the second data row down to 100
#    (matching the first data row) instead of 100000.
#  - Leave the cursor/selection on D8, matching the state the file was saved in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D header text -------------------------------------------------
$ws.Range("D1").Value = "Surface tension coefficient (sigma) (Causes stability issues if turned up above 0)"

# --- Column D values & number format (sigma values reset to 0) -----------
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("D3").NumberFormat = "0"
$ws.Range("D4").NumberFormat = "0"

# --- Widen column D to fit the new, longer header text --------------------
$ws.Columns("D").ColumnWidth = 70.5

# --- Fix x size / y size on the second data row (row 4) -------------------
$ws.Range("N4").Value = 100
$ws.Range("O4").Value = 100

# --- Restore cursor/selection position -------------------------------------
$ws.Range("D8").Select()
